# chore: publish IG 1.0.1
# - rename the "Include ..." worksheet tab to "Include #0"
# - bump Version to 1.0.1
# - update Contact
# - insert a new "Jurisdiction" row right after "Contact"

$wb = $excel.ActiveWorkbook

# Rename the second sheet's tab
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include #0"

$ws = $wb.Worksheets.Item(1)

# Bump the Version value
$ws.Range("B3").Value = "1.0.1"

# Update the Contact value
$ws.Range("B10").Value = "MedCom (http://www.medcom.dk)"

# Insert a new row after the Contact row (row 11) by first extending the
# table with one extra row (copying the formatting of the last existing
# data row so the new row matches the rest of the table), then shifting
# the values of rows 11-14 down into 12-15, leaving row 11 free for the
# new "Jurisdiction" entry.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

for ($r = 14; $r -ge 11; $r--) {
    $destRow = $r + 1
    $ws.Range("A$destRow").Value = $ws.Range("A$r").Value()
    $ws.Range("B$destRow").Value = $ws.Range("B$r").Value()
}

$ws.Range("A11").Value = "Jurisdiction"

# Give the Jurisdiction row an empty (but present) text value, matching
# the blank "value" cells used elsewhere in the workbook, by copying the
# already-blank-text cell from the Include sheet.
$wsInclude.Range("A3").Copy()
$ws.Range("B11").PasteSpecial(-4163)
